$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right after "总计" (i.e. before the
#    existing "2022-Q1" sheet). We copy an existing quarter sheet that
#    already has the identical A1:H5 layout/styling (2021-Q3, 4 data rows)
#    so the new sheet inherits matching formatting, then overwrite its
#    cell contents with the 2022-Q4 fund-holding data.
# ---------------------------------------------------------------------------
$srcTemplate = $wb.Worksheets.Item("2021-Q3")
$beforeSheet = $wb.Worksheets.Item("2022-Q1")
$srcTemplate.Copy($beforeSheet)

$q4 = $wb.Worksheets.Item("2021-Q3 (2)")
$q4.Name = "2022-Q4"

# Header row (D1 differs between templates: "基金金额" -> "基金规模")
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Row 2
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'013340"
$q4.Range("C2").Value = "创金合信芯片产业股票C"
$q4.Range("D2").Value = "'0.94"
$q4.Range("E2").Value = "'92.41"
$q4.Range("F2").Value = "'5.15"
$q4.Range("G2").Value = "'0.0484"
$q4.Range("H2").Value = 4

# Row 3
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'013339"
$q4.Range("C3").Value = "创金合信芯片产业股票A"
$q4.Range("D3").Value = "'0.92"
$q4.Range("E3").Value = "'92.41"
$q4.Range("F3").Value = "'5.15"
$q4.Range("G3").Value = "'0.0474"
$q4.Range("H3").Value = 4

# Row 4
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'519172"
$q4.Range("C4").Value = "浦银安盛睿智精选灵活配置混合A"
$q4.Range("D4").Value = "'0.21"
$q4.Range("E4").Value = "'89.04"
$q4.Range("F4").Value = "'2.77"
$q4.Range("G4").Value = "'0.0058"
$q4.Range("H4").Value = 9

# Row 5
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'519173"
$q4.Range("C5").Value = "浦银安盛睿智精选灵活配置混合C"
$q4.Range("D5").Value = "'0.13"
$q4.Range("E5").Value = "'89.04"
$q4.Range("F5").Value = "'2.77"
$q4.Range("G5").Value = "'0.0036"
$q4.Range("H5").Value = 9

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q4 and
#    renumber the leading index column (A) for every row that shifted down.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.11

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
